{"js": "// Revision 1 de la profesora\n// 1) \"Puntos positivos\" -> \"Puntos positivos.\"\n// 2) Move the _GoBack bookmark from that paragraph down to the blank\n//    paragraph right before \"Puntos negativos\".\n// 3) \"Ofrecer una gu\u00eda t\u00e9cnica detallada\" -> \"...gu\u00eda t\u00e9cnica Explicada \"\n// 4) \"Puntos negativos\" -> \"Puntos en contra \"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Append \".\" right after \"Puntos positivos\" (bold run, same formatting).\nconst positivos = paragraphs.items[0];\nconst periodRun = positivos.insertText(\".\", Word.InsertLocation.end);\nperiodRun.font.bold = true;\nawait context.sync();\n\n// 2) Relocate the _GoBack bookmark to the empty paragraph that separates\n// the \"positivos\" bullets from \"Puntos negativos\" (the 6th paragraph,\n// index 5).\ncontext.document.deleteBookmark(\"_GoBack\");\nconst blankParagraph = paragraphs.items[5];\nblankParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Replace \"detallada\" with \"Explicada \" in the gu\u00eda t\u00e9cnica bullet.\nconst detalladaResults = body.search(\"detallada\", { matchCase: true, matchWholeWord: true });\ndetalladaResults.load(\"items\");\nawait context.sync();\ndetalladaResults.items[0].insertText(\"Explicada \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Replace \"negativos\" with \"en contra \" in the \"Puntos negativos\" header.\nconst negativosResults = body.search(\"negativos\", { matchCase: true, matchWholeWord: true });\nnegativosResults.load(\"items\");\nawait context.sync();\nnegativosResults.items[0].insertText(\"en contra \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Revision 1 de la profesora\n# 1) \"Puntos positivos\" -> \"Puntos positivos.\"\n# 2) Move the _GoBack bookmark from that paragraph down to the blank\n#    paragraph right before \"Puntos negativos\".\n# 3) \"Ofrecer una gu\u00eda t\u00e9cnica detallada\" -> \"...gu\u00eda t\u00e9cnica Explicada \"\n# 4) \"Puntos negativos\" -> \"Puntos en contra \"\n\n$d = $word.ActiveDocument\n\n# 1) Append \".\" right after \"Puntos positivos\" (inherits the bold run format).\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.InsertAfter(\".\")\n\n# 2) Relocate the _GoBack bookmark to the blank paragraph (6th paragraph)\n# that separates the \"positivos\" bullets from \"Puntos negativos\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$p6 = $d.Paragraphs.Item(6)\n$d.Bookmarks.Add(\"_GoBack\", $p6.Range)\n\n# 3) Replace \"detallada\" with \"Explicada \" in the gu\u00eda t\u00e9cnica bullet.\n$find1 = $d.Content\n$find1.Find.Execute(\"detallada\", $false, $true, $false, $false, $false, $true, 1, $false, \"Explicada \", 2)\n\n# 4) Replace \"negativos\" with \"en contra \" in the \"Puntos negativos\" header.\n$find2 = $d.Content\n$find2.Find.Execute(\"negativos\", $false, $true, $false, $false, $false, $true, 1, $false, \"en contra \", 2)\n"}
